# Updated test data for DC, TripCurrent, Voltdrop, BatteryStandby
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")  # the "Add Panels" sheet (first tab) holds the edited test data

# Row 8: update panel / accessory data
$ws.Range("A8").Value = "FC64-2"
$ws.Range("C8").Value = "CPU 801"
$ws.Range("B4").Value = "NGC-1928/T961 OR TC-71696"
$ws.Range("I8").Value = "Generic Printer"
$ws.Range("K8").Value = "Printer 1"

$ws.Range("F8").Value = 0.223
$ws.Range("G8").Value = 0.415
$ws.Range("L8").Value = "'0.000"
$ws.Range("M8").Value = "'0.000"

# Row 3: add the test method / user story name
$ws.Range("B3").Value = "verifyBatteryStandbyAndAlarmLoadOnAdditionAndDeletionOfAccessories"

# Update the active selection to match the saved view state
[void]$ws.Range("B6").Select()
